$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 36014684-8719-4fb6-99ad-182db142a162.md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-21 04:54:04"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Row 2 -> 0674f237-ad4c-4385-9969-71800a4ae1ee.md : only status changes
$wsZhCn.Range("C2").Value = "Ready for handoff"
# Row 3 -> 36014684-8719-4fb6-99ad-182db142a162.md : status, handoff datetime, error detail
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-21 04:53:58"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55d7a62edcebe2709455e557e795de891c129e43/e2e/36014684-8719-4fb6-99ad-182db142a162.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6bb483946efa0ae93b663b0a64f3e734d50ceca/e2e/36014684-8719-4fb6-99ad-182db142a162.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.14

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Row 2 -> 0674f237-ad4c-4385-9969-71800a4ae1ee.md : unchanged
# Row 3 -> 36014684-8719-4fb6-99ad-182db142a162.md : status, handoff datetime, error detail
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-21 04:54:04"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55d7a62edcebe2709455e557e795de891c129e43/e2e/36014684-8719-4fb6-99ad-182db142a162.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6bb483946efa0ae93b663b0a64f3e734d50ceca/e2e/36014684-8719-4fb6-99ad-182db142a162.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.14
